$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ftests")

# --- Row 45: fm40 ---------------------------------------------------------
# Pick up the formatting already used on the neighbouring rows: columns
# B, C, H, I copy the look of the row above (row 44); columns D and G copy
# the look used a few rows up (row 29), which is the existing row with the
# same "string deductible code in D / right-aligned number in G" shape.
$ws.Range("B44").Copy()
$ws.Range("B45").PasteSpecial(-4122)
$ws.Range("C44").Copy()
$ws.Range("C45").PasteSpecial(-4122)
$ws.Range("H44").Copy()
$ws.Range("H45").PasteSpecial(-4122)
$ws.Range("I44").Copy()
$ws.Range("I45").PasteSpecial(-4122)
$ws.Range("D29").Copy()
$ws.Range("D45").PasteSpecial(-4122)
$ws.Range("G29").Copy()
$ws.Range("G45").PasteSpecial(-4122)

$ws.Cells.Item(45, 2).Value = "fm40"
$ws.Cells.Item(45, 3).Value = "OED spec example 5 - multiple policy layers"
$ws.Cells.Item(45, 4).Value = "All"
$ws.Cells.Item(45, 7).Value = 2
$ws.Cells.Item(45, 8).Value = "in progress"
$ws.Cells.Item(45, 9).Value = "in progress"

# --- Row 46: fm41 ----------------------------------------------------------
$ws.Range("B44").Copy()
$ws.Range("B46").PasteSpecial(-4122)
$ws.Range("C44").Copy()
$ws.Range("C46").PasteSpecial(-4122)
$ws.Range("H44").Copy()
$ws.Range("H46").PasteSpecial(-4122)
$ws.Range("I44").Copy()
$ws.Range("I46").PasteSpecial(-4122)
$ws.Range("D29").Copy()
$ws.Range("D46").PasteSpecial(-4122)
$ws.Range("G29").Copy()
$ws.Range("G46").PasteSpecial(-4122)

$ws.Cells.Item(46, 2).Value = "fm41"
$ws.Cells.Item(46, 3).Value = "A single special condition on locations 1 and 2, based on OED spec example 3"
$ws.Cells.Item(46, 4).Value = "All"
$ws.Cells.Item(46, 7).Value = 1
$ws.Cells.Item(46, 8).Value = "in progress"
$ws.Cells.Item(46, 9).Value = "in progress"

$excel.CutCopyMode = $false

# Match the updated selection left behind by the edit (new last row).
$ws.Activate()
$ws.Range("C46").Select() | Out-Null
